$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the two empty trailing rows (former rows 26 and 28). Delete from
#    the bottom up so row numbers of not-yet-deleted rows stay stable. This
#    also shifts the old row 27 ("8. Ссылки и документация") up to row 26.
# ---------------------------------------------------------------------------
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# ---------------------------------------------------------------------------
# 2. Give every "field label" cell in column A a light-blue fill (matching
#    the section-header fill). Mint the new format on A2 then fan it out to
#    the rest one destination block at a time, always copying from a single
#    source cell so Excel's paste never "tiles" past the destination block.
# ---------------------------------------------------------------------------
$ws.Range("A2").Interior.Color = 15652797
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A6:A10").PasteSpecial(-4122)
$ws.Range("A12:A14").PasteSpecial(-4122)
$ws.Range("A16:A17").PasteSpecial(-4122)
$ws.Range("A19:A21").PasteSpecial(-4122)
$ws.Range("A23:A24").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. B2 loses its fill (it now matches B3's plain bordered/wrapped style).
# ---------------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Add a blank, filled B cell alongside every section-header row (B1, B5,
#    B11, B15, B18, B22): Times New Roman 10pt font + the same light-blue
#    fill, no border. Mint it on B1, then copy out to the rest.
# ---------------------------------------------------------------------------
$ws.Range("B1").Font.Name = "Times New Roman"
$ws.Range("B1").Font.Size = 10
$ws.Range("B1").Interior.Color = 15652797
$ws.Range("B1").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B22").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5. Update the sheet view: scroll position and active selection.
# ---------------------------------------------------------------------------
$ws.Range("A13").Select()
$ws.Range("H9").Select()
